# Populate the new "bewerken" (edit) use-case sub-blocks and update the
# "Project bewerken" block with its final content (rows 69-118).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values are written in the exact order the original author entered them so that
# newly created shared strings land at the same index as in the target file.

$ws.Range("D70").Value = @'
use case 'Projecten toevoegen' of 'Projecten beheren' wordt uitgevoerd tot en met referentie naar deze use case
'@

$ws.Range("D108").Value = @'
Subscribertype werd gekozen voor bepaald project
'@

$ws.Range("D89").Value = @'
1. Systeem toont invoerveld om foto-URL toe te voegen
2. Gebruiker voert fotolocatie in
3. Systeem valideert invoer
4. Systeem geeft controle terug aan parent- use case
'@
$ws.Rows.Item(89).RowHeight = 60

$ws.Range("D71").Value = @'
1. Systeem laad mogelijke categorieen
2. Gebruiker kieset categorie
3. Systeem geeft controle terug aan parent- use case
'@
$ws.Rows.Item(71).RowHeight = 45

$ws.Range("D80").Value = @'
1. Systeem toont invoerveld om document-URL toe te voegen
2. Gebruiker voert documentlocatie in
3. Systeem valideert invoer
4. Systeem geeft controle terug aan parent- use case
'@
$ws.Rows.Item(80).RowHeight = 60

$ws.Range("D82").Value = @'
[invoervalidatie mislukt]
1. Systeem toont melding "URL is niet correct"
2. Systeem maakt invoer leeg en geeft controle terug aan de gebruiker
'@
$ws.Rows.Item(82).RowHeight = 45.75

$ws.Range("D81").Value = @'
Document toegevoegd/bewerkt
'@

$ws.Range("D72").Value = @'
Categorie toegevoegd/bewerkt
'@

$ws.Range("D90").Value = @'
foto toegevoegd/bewerkt
'@

$ws.Range("D99").Value = @'
video toegevoegd/bewerkt
'@

$ws.Range("D98").Value = @'
1. Systeem toont invoerveld om video-URL toe te voegen via youtube
2. Gebruiker voert videolocatie in
3. Systeem valideert invoer
4. Systeem geeft controle terug aan parent- use case
'@
$ws.Rows.Item(98).RowHeight = 60

$ws.Range("D107").Value = @'
1. Systeem toont invoerveld om subscribertype toe te voegen
2. Gebruiker kiest een type (Health/IT/…)
3. Systeem geeft controle terug aan parent- use case
'@
$ws.Rows.Item(107).RowHeight = 45

$ws.Range("D113").Value = @'
Bestaand project aanpassen
'@

$ws.Range("D114").Value = @'
Bestaande PXL student die het project heeft aangemaakt / docent die is ingeschreven voor het project
'@

$ws.Range("D116").Value = @'
1. Gebruiker druk op een knop op een bestaand project aan te passen
2. Gebruiker voert een titel, beschrijving, begin- en einddatum in
3. Gebruiker kiest een categorie (ref. use case 'Categorie kiezen')
4. Gebruiker voegt eventueel documenten toe (ref. use case 'Documenten toevoegen')
5. Gebruiker voegt eventueel foto(s) toe (ref. use case 'Foto toevoegen')
6. Gebruiker voegt eventueel video(s) toe (ref. use case 'Video toevoegen')
7. Gebruiker voert aantal subscribers in met hun type (ref. use case 'subscribertypes toevoegen')
'@
$ws.Rows.Item(116).RowHeight = 105

$ws.Range("D117").Value = @'
Er werd een bestaand project aangepast
'@

$ws.Range("D118").Value = @'
[Gebruiker annuleert aanpassing]
1. Gebruiker klikt op knop 'annuleren'
2. Systeem laat projectenlijst zien (ref. use case 'Projectenlijst tonen')
'@
$ws.Rows.Item(118).RowHeight = 45.75

$ws.Range("D69").Value = @'
Bestaande PXL student/docent
'@

$ws.Range("D79").Value = @'
use case 'Projecten toevoegen' of 'Projecten beheren' wordt uitgevoerd tot en met referentie naar deze use case
'@

$ws.Range("D88").Value = @'
use case 'Projecten toevoegen' of 'Projecten beheren' wordt uitgevoerd tot en met referentie naar deze use case
'@

$ws.Range("D97").Value = @'
use case 'Projecten toevoegen' of 'Projecten beheren' wordt uitgevoerd tot en met referentie naar deze use case
'@

$ws.Range("D106").Value = @'
use case 'Projecten toevoegen' of 'Projecten beheren' wordt uitgevoerd tot en met referentie naar deze use case
'@

$ws.Range("D78").Value = @'
Bestaande PXL student/docent
'@

$ws.Range("D87").Value = @'
Bestaande PXL student/docent
'@

$ws.Range("D96").Value = @'
Bestaande PXL student/docent
'@

$ws.Range("D105").Value = @'
Bestaande PXL student/docent
'@

$ws.Range("D91").Value = @'
[invoervalidatie mislukt]
1. Systeem toont melding "URL is niet correct"
2. Systeem maakt invoer leeg en geeft controle terug aan de gebruiker
'@
$ws.Rows.Item(91).RowHeight = 45.75

$ws.Range("D100").Value = @'
[invoervalidatie mislukt]
1. Systeem toont melding "URL is niet correct"
2. Systeem maakt invoer leeg en geeft controle terug aan de gebruiker
'@
$ws.Rows.Item(100).RowHeight = 45.75

$ws.Range("D115").Value = @'
use case 'projectenlijst tonen' werd successvol afgerond.
'@

# Update the window scroll position / active selection to match the author's
# final view (scrolled further down the sheet, with D121 as the active cell).
try {
    $excel.ActiveWindow.ScrollRow = 104
    $excel.ActiveWindow.ScrollColumn = 2
} catch {
    # Scroll position isn't always settable in every host - ignore failures.
}
$ws.Range("D121").Select()
